# Populate "TestData" style worksheet with license-plate lookup table.
# The write order below matters: it reproduces the shared-string table
# order seen in the target workbook (first-use order), which mirrors how
# the data was actually entered in Excel: plates down column A first,
# then Known Make / Known Colour filled row by row, then the header row,
# then the Retrived make / Retrived colour columns (E/F) reusing values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 2 reuses the pre-existing "G607 DAN", then Known Make ---
$ws.Cells.Item(2,1).Value = "G607 DAN"
$ws.Cells.Item(2,2).Value = "MAZDA"

# --- rest of column A (license plates) ---
$ws.Cells.Item(3,1).Value = "FV65 ZLZ"
$ws.Cells.Item(4,1).Value = "Y993 PAX"
$ws.Cells.Item(5,1).Value = "Y946 CJW"
$ws.Cells.Item(6,1).Value = "Y736 AKK"
$ws.Cells.Item(7,1).Value = "T62 EVC"

# --- Known Make / Known Colour, row by row ---
$ws.Cells.Item(3,2).Value = "KIA"
$ws.Cells.Item(3,3).Value = "WHITE"
$ws.Cells.Item(4,2).Value = "HONDA"
$ws.Cells.Item(4,3).Value = "BLACK"
$ws.Cells.Item(5,2).Value = "SUZUKI"
$ws.Cells.Item(5,3).Value = "SILVER"
$ws.Cells.Item(6,2).Value = "VAUXHALL"
$ws.Cells.Item(6,3).Value = "GREEN"
$ws.Cells.Item(2,3).Value = "BLUE"

# --- Header row ---
$ws.Cells.Item(1,1).Value = "License Plate"
$ws.Cells.Item(1,2).Value = "Known Make"
$ws.Cells.Item(1,3).Value = "Known Colour"
$ws.Cells.Item(1,5).Value = "Retrived make"
$ws.Cells.Item(1,6).Value = "Retrived colour"

# --- Row 7 Known Make/Colour (introduces PEUGEOT) ---
$ws.Cells.Item(7,2).Value = "PEUGEOT"
$ws.Cells.Item(7,3).Value = "BLUE"

# --- Retrived make / Retrived colour columns (E/F), reusing known values ---
$ws.Cells.Item(2,5).Value = "MAZDA"
$ws.Cells.Item(2,6).Value = "BLUE"
$ws.Cells.Item(3,5).Value = "KIA"
$ws.Cells.Item(3,6).Value = "WHITE"
$ws.Cells.Item(4,5).Value = "HONDA"
$ws.Cells.Item(4,6).Value = "BLACK"
$ws.Cells.Item(5,5).Value = "SUZUKI"
$ws.Cells.Item(5,6).Value = "SILVER"
$ws.Cells.Item(6,5).Value = "VAUXHALL"
$ws.Cells.Item(6,6).Value = "GREEN"
$ws.Cells.Item(7,5).Value = "PEUGEOT"
$ws.Cells.Item(7,6).Value = "BLUE"

# --- Column widths (best-effort: runtime quantizes to its own pixel grid) ---
$ws.Columns.Item(1).ColumnWidth = 14.14
$ws.Columns.Item(2).ColumnWidth = 14.71
$ws.Columns.Item(3).ColumnWidth = 14.43
$ws.Columns.Item(5).ColumnWidth = 15.29
$ws.Columns.Item(6).ColumnWidth = 14.14

# --- Selection matches the authored state: rectangle E2:F9 active at E9 ---
$ws.Range("E2:F9").Select()
